# "Bit fix of Presentation.pptx - year and name of the game"
#
# 1) Slide 1 (title slide): the game name changes from the single run
#    "Five-Cards Draw" into two runs "Five Card " + "Draw" (same
#    formatting carried over to both runs).
# 2) Slide 13 (credits slide): the year textbox changes from "2014" to
#    "2015".

$p = $ppt.ActivePresentation

# --- Slide 1: fix the game's name -----------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(2)          # "Подзаглавие 2" subTitle placeholder
$fullRange = $titleShape.TextFrame.TextRange
# Original text: "Five-Cards Draw" (15 chars).
# Replace the first 11 characters ("Five-Cards ") with "Five Card " so
# that the remaining text is exactly "Draw" in its own run, matching
# the two-run split produced by the authoring edit.
$firstPart = $fullRange.Characters(1, 11)
$firstPart.Text = "Five Card "

# --- Slide 13: fix the year ------------------------------------------
$slide13 = $p.Slides.Item(13)
$yearShape = $slide13.Shapes.Item(4)          # "Текстово поле 6" textbox with "2014"
$yearShape.TextFrame.TextRange.Text = "2015"
